# Pb208_RMF.xlsx edit: fix Dirac-equation solver input guesses so that
# O16, Ca40, and Pb208 converge without needing 0 iterations.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Sheet1
$ws2 = $wb.Worksheets.Item(2)   # Sheet2
$ws3 = $wb.Worksheets.Item(3)   # Sheet3 (Pb208, kappa=+1/2 family)
$ws4 = $wb.Worksheets.Item(4)   # Sheet4 (Pb208, kappa=-1/2 family)

# --- Sheet1: N_STEPS 200 -> 375 ---
$ws1.Range("G2").Value = 375

# --- Sheet3 cell updates ---
$ws3.Range("C2").Formula = "=878"
$ws3.Range("E2").Value = 5
$ws3.Range("E3").Value = 5
$ws3.Range("E4").Value = 2
$ws3.Range("E5").Value = 6
$ws3.Range("A6").Value = 4
$ws3.Range("B6").Value = 2
$ws3.Range("D6").Value = "1D3/2"
$ws3.Range("E6").Value = 3
$ws3.Range("A7").Value = 2
$ws3.Range("B7").Value = -1
$ws3.Range("D7").Value = "2S1/2"
$ws3.Range("E7").Value = 6
$ws3.Range("E8").Value = 6
$ws3.Range("A9").Value = 6
$ws3.Range("B9").Value = 3
$ws3.Range("D9").Value = "1F5/2"
$ws3.Range("E9").Value = 4
$ws3.Range("A10").Value = 4
$ws3.Range("B10").Value = -2
$ws3.Range("D10").Value = "2P3/2"
$ws3.Range("E10").Value = 7
$ws3.Range("E11").Value = 2
$ws3.Range("E12").Value = 6
$ws3.Range("E13").Value = 4
$ws3.Range("E14").Value = 7
$ws3.Range("A15").Value = 12
$ws3.Range("B15").Value = -6
$ws3.Range("D15").Value = "1H11/2"
$ws3.Range("E15").Value = 6
$ws3.Range("A16").Value = 4
$ws3.Range("B16").Value = 2
$ws3.Range("D16").Value = "2D3/2"
$ws3.Range("E16").Value = 2
$ws3.Range("E17").Value = 2

# --- Sheet4 cell updates ---
$ws4.Range("C2").Formula = "=868"
$ws4.Range("E2").Value = 5
$ws4.Range("E3").Value = 5
$ws4.Range("E4").Value = 2
$ws4.Range("E5").Value = 6
$ws4.Range("A6").Value = 4
$ws4.Range("B6").Value = 2
$ws4.Range("D6").Value = "1D3/2"
$ws4.Range("E6").Value = 4
$ws4.Range("A7").Value = 2
$ws4.Range("B7").Value = -1
$ws4.Range("D7").Value = "2S1/2"
$ws4.Range("E7").Value = 6
$ws4.Range("E8").Value = 6
$ws4.Range("A9").Value = 6
$ws4.Range("B9").Value = 3
$ws4.Range("D9").Value = "1F5/2"
$ws4.Range("E9").Value = 4
$ws4.Range("A10").Value = 4
$ws4.Range("B10").Value = -2
$ws4.Range("D10").Value = "2P3/2"
$ws4.Range("E10").Value = 7
$ws4.Range("E11").Value = 2
$ws4.Range("E12").Value = 6
$ws4.Range("E13").Value = 4
$ws4.Range("A14").Value = 6
$ws4.Range("B14").Value = -3
$ws4.Range("D14").Value = "2D5/2"
$ws4.Range("E14").Value = 3
$ws4.Range("A15").Value = 12
$ws4.Range("B15").Value = -6
$ws4.Range("D15").Value = "1H11/2"
$ws4.Range("E15").Value = 6
$ws4.Range("E16").Value = 2
$ws4.Range("A17").Value = 2
$ws4.Range("B17").Value = -1
$ws4.Range("D17").Value = "3S1/2"
$ws4.Range("E17").Value = 8
$ws4.Range("A18").Value = 10
$ws4.Range("B18").Value = 5
$ws4.Range("D18").Value = "1H9/2"
$ws4.Range("E18").Value = 5
$ws4.Range("E19").Value = 4
$ws4.Range("A20").Value = 14
$ws4.Range("B20").Value = -7
$ws4.Range("D20").Value = "1I13/2"
$ws4.Range("E21").Value = 3
$ws4.Range("A22").Value = 4
$ws4.Range("B22").Value = -2
$ws4.Range("D22").Value = "3P3/2"
$ws4.Range("E22").Value = 8
$ws4.Range("A23").Value = 2
$ws4.Range("B23").Value = 1
$ws4.Range("D23").Value = "3P1/2"
$ws4.Range("E23").Value = 1

# --- Selections ---
# Set selections on non-active sheets first, then finish on Sheet3 so it
# remains the active tab (matches activeTab=2 / tabSelected on Sheet3).
$ws2.Select()
$ws2.Range("F2").Select() | Out-Null

$ws4.Select()
$ws4.Range("E36").Select() | Out-Null

$ws3.Select()
$ws3.Range("C2").Select() | Out-Null
